$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 21
$ws.Range("H21").Value = 21249.5
$ws.Range("I21").Value = 10000
$ws.Range("J21").Value = 24999.334
$ws.Range("K21").Value = 10000
$ws.Range("L21").Value = 24999.334
$ws.Range("M21").Value = -9532
$ws.Range("N21").Value = -25935.334
# row 23
$ws.Range("H23").Value = 21249.5
$ws.Range("I23").Value = 10000
$ws.Range("J23").Value = 24999.334
$ws.Range("K23").Value = 10000
$ws.Range("L23").Value = 24999.334
$ws.Range("M23").Value = -9766
$ws.Range("N23").Value = -25467.334
# row 29
$ws.Range("H29").Value = 200
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = ""
# row 31
$ws.Range("H31").Value = 3000
$ws.Range("J31").Value = 3000
$ws.Range("L31").Value = 9000
$ws.Range("N31").Value = -9460
# row 38
$ws.Range("H38").Value = 422.8125
$ws.Range("I38").Value = 319.69232
$ws.Range("J38").Value = 869.6667
$ws.Range("K38").Value = 959.07696
$ws.Range("L38").Value = 2609.0001
$ws.Range("M38").Value = -587.07696
$ws.Range("N38").Value = -3353.0001
# row 58
$ws.Range("H58").Value = 1365.625
$ws.Range("J58").Value = 2000
$ws.Range("L58").Value = 6000
$ws.Range("N58").Value = -6300
# row 87
$ws.Range("H87").Value = 23306.062
$ws.Range("J87").Value = 23306.062
$ws.Range("L87").Value = 23306.062
$ws.Range("N87").Value = -25802.062
# row 90
$ws.Range("H90").Value = 23306.062
$ws.Range("J90").Value = 23306.062
$ws.Range("L90").Value = 69918.186
$ws.Range("N90").Value = -82398.186
# row 129
$ws.Range("H129").Value = 1339.381
$ws.Range("J129").Value = 1937.3636
$ws.Range("L129").Value = 5812.0908
$ws.Range("N129").Value = -15812.0908
# row 137
$ws.Range("H137").Value = 1334.4783
$ws.Range("I137").Value = 1353.0769
$ws.Range("J137").Value = 1310.3
$ws.Range("K137").Value = 4059.2307
$ws.Range("L137").Value = 3930.9
$ws.Range("M137").Value = -1509.2307
$ws.Range("N137").Value = -9030.9
# row 138
$ws.Range("H138").Value = 2712.3655
$ws.Range("I138").Value = 3517
$ws.Range("J138").Value = 2496.4878
$ws.Range("K138").Value = 10551
$ws.Range("L138").Value = 7489.4634
$ws.Range("M138").Value = -5411
$ws.Range("N138").Value = -17769.4634

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 17
$ws.Range("H17").Value = 10000
$ws.Range("I17").Value = 10000
$ws.Range("K17").Value = 10000
$ws.Range("M17").Value = -9827
# row 40
$ws.Range("H40").Value = 5866.6665
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 5866.6665
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 5866.6665
$ws.Range("M40").Value = ""
$ws.Range("N40").Value = -6218.6665
# row 61
$ws.Range("H61").Value = 11496722
$ws.Range("I61").Value = 33334636
$ws.Range("K61").Value = 33334636
$ws.Range("M61").Value = -33334424
# row 74
$ws.Range("H74").Value = 1237.9546
$ws.Range("I74").Value = 701.8
$ws.Range("K74").Value = 701.8
$ws.Range("M74").Value = 172.2
# row 77
$ws.Range("H77").Value = 1237.9546
$ws.Range("I77").Value = 701.8
$ws.Range("K77").Value = 3509
$ws.Range("M77").Value = 859
# row 136
$ws.Range("H136").Value = 11496722
$ws.Range("I136").Value = 33334636
$ws.Range("K136").Value = 100003908
$ws.Range("M136").Value = -100001358

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 80
$ws.Range("H80").Value = 2264.7646
$ws.Range("I80").Value = 2763.4285
$ws.Range("J80").Value = 1915.7
$ws.Range("K80").Value = 2763.4285
$ws.Range("L80").Value = 1915.7
$ws.Range("M80").Value = -1765.4285
$ws.Range("N80").Value = -3911.7
# row 83
$ws.Range("H83").Value = 2264.7646
$ws.Range("I83").Value = 2763.4285
$ws.Range("J83").Value = 1915.7
$ws.Range("K83").Value = 13817.1425
$ws.Range("L83").Value = 9578.5
$ws.Range("M83").Value = -8825.1425
$ws.Range("N83").Value = -19562.5
# row 94
$ws.Range("H94").Value = 993.7
$ws.Range("I94").Value = 905.2857
$ws.Range("J94").Value = 1200
$ws.Range("K94").Value = 905.2857
$ws.Range("L94").Value = 1200
$ws.Range("M94").Value = -454.2857
$ws.Range("N94").Value = -2102
# row 134
$ws.Range("H134").Value = 2562.9473
$ws.Range("I134").Value = 2346.4443
$ws.Range("J134").Value = 3094.3635
$ws.Range("K134").Value = 7039.3329
$ws.Range("L134").Value = 9283.0905
$ws.Range("M134").Value = -4504.3329
$ws.Range("N134").Value = -14353.0905

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 23
$ws.Range("H23").Value = 14005.556
$ws.Range("I23").Value = 6625
$ws.Range("J23").Value = 19910
$ws.Range("K23").Value = 6625
$ws.Range("L23").Value = 19910
$ws.Range("M23").Value = -6385
$ws.Range("N23").Value = -20390
# row 27
$ws.Range("H27").Value = 14005.556
$ws.Range("I27").Value = 6625
$ws.Range("J27").Value = 19910
$ws.Range("K27").Value = 6625
$ws.Range("L27").Value = 19910
$ws.Range("M27").Value = -6433
$ws.Range("N27").Value = -20294
# row 31
$ws.Range("H31").Value = 3814.52
$ws.Range("I31").Value = 1434.8334
$ws.Range("J31").Value = 5153.0938
$ws.Range("K31").Value = 1434.8334
$ws.Range("L31").Value = 5153.0938
$ws.Range("M31").Value = -1139.8334
$ws.Range("N31").Value = -5743.0938
# row 34
$ws.Range("H34").Value = 3814.52
$ws.Range("I34").Value = 1434.8334
$ws.Range("J34").Value = 5153.0938
$ws.Range("K34").Value = 1434.8334
$ws.Range("L34").Value = 5153.0938
$ws.Range("M34").Value = -1232.8334
$ws.Range("N34").Value = -5557.0938
# row 58
$ws.Range("H58").Value = 2627.5293
$ws.Range("I58").Value = 2430.6667
$ws.Range("J58").Value = 3100
$ws.Range("K58").Value = 2430.6667
$ws.Range("L58").Value = 3100
$ws.Range("M58").Value = -2227.6667
$ws.Range("N58").Value = -3506
# row 132
$ws.Range("H132").Value = 11907284
$ws.Range("I132").Value = 1468.5
$ws.Range("J132").Value = 20836646
$ws.Range("K132").Value = 4405.5
$ws.Range("L132").Value = 62509938
$ws.Range("M132").Value = -1875.5
$ws.Range("N132").Value = -62514998
# row 134
$ws.Range("H134").Value = 968.3333
$ws.Range("I134").Value = 718.73914
$ws.Range("K134").Value = 2156.21742
$ws.Range("M134").Value = 378.7825800000001
# row 136
$ws.Range("H136").Value = 2627.5293
$ws.Range("I136").Value = 2430.6667
$ws.Range("J136").Value = 3100
$ws.Range("K136").Value = 7292.000100000001
$ws.Range("L136").Value = 9300
$ws.Range("M136").Value = -4742.000100000001
$ws.Range("N136").Value = -14400

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 68
$ws.Range("H68").Value = 1424.1609
$ws.Range("I68").Value = 810.4483
$ws.Range("J68").Value = 1731.0172
$ws.Range("K68").Value = 2431.3449
$ws.Range("L68").Value = 5193.0516
$ws.Range("M68").Value = -1620.3449
$ws.Range("N68").Value = -6815.0516
# row 71
$ws.Range("H71").Value = 1424.1609
$ws.Range("I71").Value = 810.4483
$ws.Range("J71").Value = 1731.0172
$ws.Range("K71").Value = 7294.0347
$ws.Range("L71").Value = 15579.1548
$ws.Range("M71").Value = -3238.0347
$ws.Range("N71").Value = -23691.1548
# row 107
$ws.Range("H107").Value = 1559.1971
$ws.Range("I107").Value = 319.5862
$ws.Range("J107").Value = 2415.1191
$ws.Range("K107").Value = 958.7586000000001
$ws.Range("L107").Value = 7245.3573
$ws.Range("M107").Value = 961.2413999999999
$ws.Range("N107").Value = -11085.3573
# row 110
$ws.Range("H110").Value = 12915.875
$ws.Range("I110").Value = 2109
$ws.Range("J110").Value = 19400
$ws.Range("K110").Value = 6327
$ws.Range("L110").Value = 58200
$ws.Range("M110").Value = -2237
$ws.Range("N110").Value = -66380

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 122
$ws.Range("H122").Value = 5024.625
$ws.Range("I122").Value = 4402.3335
$ws.Range("J122").Value = 5398
$ws.Range("K122").Value = 13207.0005
$ws.Range("L122").Value = 16194
$ws.Range("M122").Value = -10757.0005
$ws.Range("N122").Value = -21094
# row 126
$ws.Range("H126").Value = 1593
$ws.Range("I126").Value = 1503.4286
$ws.Range("J126").Value = 1749.75
$ws.Range("K126").Value = 4510.2858
$ws.Range("L126").Value = 5249.25
$ws.Range("M126").Value = -2040.2858
$ws.Range("N126").Value = -10189.25
# row 132
$ws.Range("H132").Value = 2089.04
$ws.Range("I132").Value = 1675.8572
$ws.Range("J132").Value = 2614.9092
$ws.Range("K132").Value = 5027.571599999999
$ws.Range("L132").Value = 7844.7276
$ws.Range("M132").Value = -2497.571599999999
$ws.Range("N132").Value = -12904.7276

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 132
$ws.Range("H132").Value = 2503.4707
$ws.Range("I132").Value = 2129.3235
$ws.Range("J132").Value = 3251.7646
$ws.Range("K132").Value = 6387.970499999999
$ws.Range("L132").Value = 9755.293799999999
$ws.Range("M132").Value = -3857.970499999999
$ws.Range("N132").Value = -14815.2938
# row 136
$ws.Range("H136").Value = 2825674.5
$ws.Range("I136").Value = 1026.4814
$ws.Range("K136").Value = 3079.4442
$ws.Range("M136").Value = -529.4441999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 14
$ws.Range("H14").Value = 5346123
$ws.Range("I14").Value = 2993
$ws.Range("K14").Value = 2993
$ws.Range("M14").Value = -2825
# row 42
$ws.Range("H42").Value = 64313.43
$ws.Range("J42").Value = 64313.43
$ws.Range("L42").Value = 64313.43
$ws.Range("N42").Value = -65069.43
# row 45
$ws.Range("H45").Value = 10695.833
$ws.Range("I45").Value = 8784.5
$ws.Range("J45").Value = 11651.5
$ws.Range("K45").Value = 8784.5
$ws.Range("L45").Value = 11651.5
$ws.Range("M45").Value = -8293.5
$ws.Range("N45").Value = -12633.5
# row 108
$ws.Range("H108").Value = 30000
$ws.Range("J108").Value = 30000
$ws.Range("L108").Value = 30000
$ws.Range("N108").Value = -37680
# row 132
$ws.Range("H132").Value = 6946652.5
$ws.Range("I132").Value = 2286.6
$ws.Range("J132").Value = 13259713
$ws.Range("K132").Value = 6859.799999999999
$ws.Range("L132").Value = 39779139
$ws.Range("M132").Value = -4329.799999999999
$ws.Range("N132").Value = -39784199
# row 136
$ws.Range("H136").Value = 2253.2295
$ws.Range("I136").Value = 2137.4546
$ws.Range("J136").Value = 2552.8823
$ws.Range("K136").Value = 6412.3638
$ws.Range("L136").Value = 7658.646900000001
$ws.Range("M136").Value = -3862.3638
$ws.Range("N136").Value = -12758.6469
